# Generate Report for Handback
$wb = $excel.ActiveWorkbook

# --- Overview sheet: status text for zh-cn / de-de columns ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("I2").Value = "7925027b-ecc4-4ad7-80b0-2392b7823ebf.md"
$wsZh.Range("J2").Value = "7925027b-ecc4-4ad7-80b0-2392b7823ebf.2a0350b9cf9f1a3d823a43ebe3722841c05dff22.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-12 03:19:57"
$zhUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/6276bc47b0fafd1404f2c9f5d69e29b334804740/e2e/7925027b-ecc4-4ad7-80b0-2392b7823ebf.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $zhUrl, [Type]::Missing, [Type]::Missing, "7925027b-ecc4-4ad7-80b0-2392b7823ebf.md")

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("I2").Value = "7925027b-ecc4-4ad7-80b0-2392b7823ebf.md"
$wsDe.Range("J2").Value = "7925027b-ecc4-4ad7-80b0-2392b7823ebf.2a0350b9cf9f1a3d823a43ebe3722841c05dff22.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-12 03:20:10"
$deUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/6276bc47b0fafd1404f2c9f5d69e29b334804740/e2e/7925027b-ecc4-4ad7-80b0-2392b7823ebf.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $deUrl, [Type]::Missing, [Type]::Missing, "7925027b-ecc4-4ad7-80b0-2392b7823ebf.md")

# --- Column width adjustments to match regenerated report layout ---
$wsOverview.Range("E1").ColumnWidth = 29.9777047293527
$wsOverview.Range("F1").ColumnWidth = 29.9777047293527

$wsZh.Range("C1").ColumnWidth = 29.9777047293527
$wsZh.Range("I1").ColumnWidth = 40
$wsZh.Range("J1").ColumnWidth = 40

$wsDe.Range("C1").ColumnWidth = 29.9777047293527
$wsDe.Range("I1").ColumnWidth = 40
$wsDe.Range("J1").ColumnWidth = 40
